$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.450.48'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '3.104.36'
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '388.55'
$ws.Range("E5").Value = '  +2.01%  '
$ws.Range("D6").Value = '103.98'
$ws.Range("E6").Value = '  +1.04%  '
$ws.Range("E7").Value = '  -1.31%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").Value = '37.18'
$ws.Range("E10").Value = '  +1.00%  '
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0860'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '3.596.56'
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("D14").Value = '18.58'
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = '3.107.31'
$ws.Range("E16").Value = '  +2.49%  '
$ws.Range("E17").Value = '  +2.30%  '
$ws.Range("D18").Value = '10.63'
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").Value = '51.553.47'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("E20").Value = '  +6.54%  '
$ws.Range("D21").Value = '12.53'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").Value = '70.32'
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("D24").Value = '266.83'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").Value = '3.18'
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").Value = '8.06'
$ws.Range("E26").Value = '  -2.45%  '
$ws.Range("D27").Value = '27.45'
$ws.Range("E27").Value = '  +4.46%  '
$ws.Range("E28").Value = '  -5.04%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -4.88%  '
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("D33").Value = '35.85'
$ws.Range("E33").Value = '  +5.24%  '
$ws.Range("D34").Value = '0.0475'
$ws.Range("E34").Value = '  +6.06%  '
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("D36").Value = '50.05'
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = '3.38'
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.290'
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("D40").Value = '129.49'
$ws.Range("E40").Value = '  +4.63%  '
$ws.Range("D41").Value = '1.86'
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.60'
$ws.Range("E42").Value = '  -2.92%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '0.116'
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").Value = '3.82'
$ws.Range("E44").Value = '  +1.33%  '
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").Value = '22.23'
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("E47").Value = '  +3.79%  '
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("D49").Value = '2.080.77'
$ws.Range("E49").Value = '  +2.25%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.933'
$ws.Range("E50").Value = '  +19.38%  '
$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").Value = '0.0331'
$ws.Range("E51").Value = '  +3.48%  '
